$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 112, shifting existing rows 112..204 down to 113..205
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new weekly record
$ws.Cells.Item(112, 1).Value = 4
$ws.Cells.Item(112, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(112, 3).Value = "Los Lagos"
$ws.Cells.Item(112, 4).Value = 44586
$ws.Cells.Item(112, 5).Value = 10
$ws.Cells.Item(112, 6).Value = 100112021
$ws.Cells.Item(112, 7).Value = "Ají"
$ws.Cells.Item(112, 8).Value = "Inferno"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 140
$ws.Cells.Item(112, 11).Value = 20000
$ws.Cells.Item(112, 12).Value = 20000
$ws.Cells.Item(112, 13).Value = 20000
$ws.Cells.Item(112, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(112, 15).Value = "Región Metropolitana"
$ws.Cells.Item(112, 16).Value = 1333
$ws.Cells.Item(112, 17).Value = 15
$ws.Cells.Item(112, 18).Value = "Hortaliza"
